$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 12, pushing the existing "Programa resumido:" block
# (and everything after it) down from rows 12-22 to rows 16-26.
$ws.Range("A12:A15").EntireRow.Insert()

# The insert leaves behind default column-styled empty cells across A:C for each new
# row; clear them so we can recreate only the cells that the target layout needs.
$ws.Range("A12:C15").Clear()

# Row 12: "Docentes responsáveis:" header - column A only, bold/top style (same as A11).
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Docentes responsáveis:"

# Rows 13-15: the three instructors, duplicated into column B (normal wrap) and
# column C (red wrap) - matching the styles used by B10/C10.
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("B13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

$ws.Range("B14").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C14").Value = "2166002 - Sandra Giacomin Schneider"

$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"

$excel.CutCopyMode = $false
